$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "43.933.30"
$ws.Cells.Item(2, 5).Value = "  +2.51%  "
# Row 3
$ws.Cells.Item(3, 4).Value = "2.346.61"
$ws.Cells.Item(3, 5).Value = "  +2.61%  "
# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.02%  "
# Row 5
$ws.Cells.Item(5, 4).Value = "313.35"
$ws.Cells.Item(5, 5).Value = "  -0.64%  "
# Row 6
$ws.Cells.Item(6, 4).Value = "'109.80"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +5.54%  "
# Row 7
$ws.Cells.Item(7, 4).Value = "0.633"
$ws.Cells.Item(7, 5).Value = "  +1.60%  "
# Row 8
$ws.Cells.Item(8, 5).Value = "  -0.01%  "
# Row 9
$ws.Cells.Item(9, 4).Value = "0.623"
$ws.Cells.Item(9, 5).Value = "  +3.67%  "
# Row 10
$ws.Cells.Item(10, 4).Value = "41.55"
$ws.Cells.Item(10, 5).Value = "  +5.58%  "
# Row 11
$ws.Cells.Item(11, 4).Value = "0.0923"
$ws.Cells.Item(11, 5).Value = "  +2.32%  "
# Row 12
$ws.Cells.Item(12, 4).Value = "8.65"
$ws.Cells.Item(12, 5).Value = "  +2.96%  "
# Row 13
$ws.Cells.Item(13, 4).Value = "1.02"
$ws.Cells.Item(13, 5).Value = "  +1.90%  "
# Row 14
$ws.Cells.Item(14, 5).Value = "  -0.79%  "
# Row 15
$ws.Cells.Item(15, 4).Value = "'15.60"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +2.65%  "
# Row 16
$ws.Cells.Item(16, 4).Value = "2.698.13"
$ws.Cells.Item(16, 5).Value = "  +2.34%  "
# Row 17
$ws.Cells.Item(17, 4).Value = "2.342.21"
$ws.Cells.Item(17, 5).Value = "  +3.00%  "
# Row 18
$ws.Cells.Item(18, 4).Value = "43.883.84"
$ws.Cells.Item(18, 5).Value = "  +2.65%  "
# Row 19
$ws.Cells.Item(19, 4).Value = "7.63"
$ws.Cells.Item(19, 5).Value = "  +2.94%  "
# Row 20
$ws.Cells.Item(20, 5).Value = "  +1.97%  "
# Row 21
$ws.Cells.Item(21, 4).Value = "13.13"
$ws.Cells.Item(21, 5).Value = "  -1.44%  "
# Row 22
$ws.Cells.Item(22, 4).Value = "74.41"
$ws.Cells.Item(22, 5).Value = "  +0.78%  "
# Row 23
$ws.Cells.Item(23, 4).Value = "3.53"
$ws.Cells.Item(23, 5).Value = "  -0.69%  "
# Row 24
$ws.Cells.Item(24, 4).Value = "269.89"
$ws.Cells.Item(24, 5).Value = "  +3.01%  "
# Row 25
$ws.Cells.Item(25, 4).Value = "2.29"
$ws.Cells.Item(25, 5).Value = "  +4.18%  "
# Row 26
$ws.Cells.Item(26, 5).Value = "  -0.16%  "
# Row 27
$ws.Cells.Item(27, 4).Value = "7.68"
$ws.Cells.Item(27, 5).Value = "  +8.91%  "
# Row 28
$ws.Cells.Item(28, 4).Value = "11.21"
$ws.Cells.Item(28, 5).Value = "  +3.20%  "
# Row 29
$ws.Cells.Item(29, 5).Value = "  -2.41%  "
# Row 30
$ws.Cells.Item(30, 4).Value = "39.19"
$ws.Cells.Item(30, 5).Value = "  +4.87%  "
# Row 31
$ws.Cells.Item(31, 4).Value = "22.73"
$ws.Cells.Item(31, 5).Value = "  +2.21%  "
# Row 32
$ws.Cells.Item(32, 4).Value = "168.39"
$ws.Cells.Item(32, 5).Value = "  +1.18%  "
# Row 33
$ws.Cells.Item(33, 4).Value = "0.0888"
$ws.Cells.Item(33, 5).Value = "  +2.03%  "
# Row 34
$ws.Cells.Item(34, 5).Value = "  +7.10%  "
# Row 35
$ws.Cells.Item(35, 5).Value = "  +1.77%  "
# Row 36
$ws.Cells.Item(36, 4).Value = "4.81"
$ws.Cells.Item(36, 5).Value = "  +5.88%  "
# Row 37
$ws.Cells.Item(37, 4).Value = "0.113"
$ws.Cells.Item(37, 5).Value = "  -0.53%  "
# Row 38
$ws.Cells.Item(38, 4).Value = "0.0368"
$ws.Cells.Item(38, 5).Value = "  +5.40%  "
# Row 39
$ws.Cells.Item(39, 4).Value = "3.84"
$ws.Cells.Item(39, 5).Value = "  +0.01%  "
# Row 40
$ws.Cells.Item(40, 4).Value = "2.88"
$ws.Cells.Item(40, 5).Value = "  +8.43%  "
# Row 41
$ws.Cells.Item(41, 4).Value = "1.72"
$ws.Cells.Item(41, 5).Value = "  +10.01%  "
# Row 42
$ws.Cells.Item(42, 4).Value = "104.58"
$ws.Cells.Item(42, 5).Value = "  +13.76%  "
# Row 43
$ws.Cells.Item(43, 4).Value = "0.239"
$ws.Cells.Item(43, 5).Value = "  +3.68%  "
# Row 44
$ws.Cells.Item(44, 2).Value = "MultiversX"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Cells.Item(44, 4).Value = "71.96"
$ws.Cells.Item(44, 5).Value = "  +3.80%  "
# Row 45
$ws.Cells.Item(45, 2).Value = "Celestia"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Cells.Item(45, 4).Value = "13.36"
$ws.Cells.Item(45, 5).Value = "  +10.09%  "
# Row 46
$ws.Cells.Item(46, 5).Value = "  -0.16%  "
# Row 47
$ws.Cells.Item(47, 4).Value = "114.65"
$ws.Cells.Item(47, 5).Value = "  +0.96%  "
# Row 48
$ws.Cells.Item(48, 4).Value = "1.669.80"
$ws.Cells.Item(48, 5).Value = "  -2.89%  "
# Row 49
$ws.Cells.Item(49, 4).Value = "77.75"
$ws.Cells.Item(49, 5).Value = "  -0.89%  "
# Row 50
$ws.Cells.Item(50, 2).Value = "FraxShare"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(50, 4).Value = "9.02"
$ws.Cells.Item(50, 5).Value = "  +2.95%  "
# Row 51
$ws.Cells.Item(51, 2).Value = "MinaProtocolToken"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina"
$ws.Cells.Item(51, 4).Value = "1.57"
$ws.Cells.Item(51, 5).Value = "  +12.19%  "

Write-Host "Applied all cell updates"
